$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total Abertura (H) / Total Clique (I) are stored as text in this sheet, so we
# briefly mark the cell as Text before writing the digits (otherwise Excel's
# normal type-inference would store the value as a number), then restore the
# default "Normal" style so no visible formatting/style change is introduced.
$ws.Range("H54").NumberFormat = "@"
$ws.Range("H54").Value = "3620"
$ws.Range("H54").Style = "Normal"

$ws.Range("H74").NumberFormat = "@"
$ws.Range("H74").Value = "3518"
$ws.Range("H74").Style = "Normal"

$ws.Range("H81").NumberFormat = "@"
$ws.Range("H81").Value = "3155"
$ws.Range("H81").Style = "Normal"

$ws.Range("H87").NumberFormat = "@"
$ws.Range("H87").Value = "3330"
$ws.Range("H87").Style = "Normal"

$ws.Range("H97").NumberFormat = "@"
$ws.Range("H97").Value = "3587"
$ws.Range("H97").Style = "Normal"

$ws.Range("H98").NumberFormat = "@"
$ws.Range("H98").Value = "243"
$ws.Range("H98").Style = "Normal"

$ws.Range("H104").NumberFormat = "@"
$ws.Range("H104").Value = "2934"
$ws.Range("H104").Style = "Normal"

$ws.Range("I104").NumberFormat = "@"
$ws.Range("I104").Value = "411"
$ws.Range("I104").Style = "Normal"

$ws.Range("H118").NumberFormat = "@"
$ws.Range("H118").Value = "366"
$ws.Range("H118").Style = "Normal"

$ws.Range("H125").NumberFormat = "@"
$ws.Range("H125").Value = "3040"
$ws.Range("H125").Style = "Normal"

$ws.Range("H126").NumberFormat = "@"
$ws.Range("H126").Value = "2159"
$ws.Range("H126").Style = "Normal"

$ws.Range("H127").NumberFormat = "@"
$ws.Range("H127").Value = "359"
$ws.Range("H127").Style = "Normal"

$ws.Range("H128").NumberFormat = "@"
$ws.Range("H128").Value = "91"
$ws.Range("H128").Style = "Normal"

$ws.Range("H133").NumberFormat = "@"
$ws.Range("H133").Value = "2697"
$ws.Range("H133").Style = "Normal"

$ws.Range("I133").NumberFormat = "@"
$ws.Range("I133").Value = "279"
$ws.Range("I133").Style = "Normal"

$ws.Range("H136").NumberFormat = "@"
$ws.Range("H136").Value = "337"
$ws.Range("H136").Style = "Normal"

$ws.Range("H137").NumberFormat = "@"
$ws.Range("H137").Value = "356"
$ws.Range("H137").Style = "Normal"

$ws.Range("H138").NumberFormat = "@"
$ws.Range("H138").Value = "301"
$ws.Range("H138").Style = "Normal"

$ws.Range("G139").Value = 901

$ws.Range("H139").NumberFormat = "@"
$ws.Range("H139").Value = "295"
$ws.Range("H139").Style = "Normal"

$ws.Range("H140").NumberFormat = "@"
$ws.Range("H140").Value = "82"
$ws.Range("H140").Style = "Normal"

$ws.Range("H141").NumberFormat = "@"
$ws.Range("H141").Value = "338"
$ws.Range("H141").Style = "Normal"

$ws.Range("G142").Value = 50641

$ws.Range("H142").NumberFormat = "@"
$ws.Range("H142").Value = "1894"
$ws.Range("H142").Style = "Normal"

$ws.Range("I142").NumberFormat = "@"
$ws.Range("I142").Value = "24"
$ws.Range("I142").Style = "Normal"

$ws.Range("H143").NumberFormat = "@"
$ws.Range("H143").Value = "285"
$ws.Range("H143").Style = "Normal"

$ws.Range("I143").NumberFormat = "@"
$ws.Range("I143").Value = "2"
$ws.Range("I143").Style = "Normal"

$ws.Range("H144").NumberFormat = "@"
$ws.Range("H144").Value = "276"
$ws.Range("H144").Style = "Normal"
